$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order of first-use for brand-new shared strings matters:
#   "MuSCs" must become the first newly-added shared string (index 21)
#   "ECs"   must become the last newly-added shared string (index 25)
# So set A3 ("MuSCs") before anything else new, and set the "ECs"
# target-cluster cells (D2/D3/D4) only at the very end.

# --- Row 2: Inflammatory-Mac | Gdf3 | Tdgf1 | (target set later) ---
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Gdf3"
$ws.Range("C2").Value = "Tdgf1"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 15.042283
$ws.Range("H2").Value = 45.126849
$ws.Range("I2").Value = 0.4622287068730439
$ws.Range("J2").Value = 0.4626185066002986
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.08586166666666667
$ws.Range("N2").Value = 0.257585
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1.291555488851667
$ws.Range("R2").Value = 11.623999399665
$ws.Range("S2").Value = 0.4622287068730439
$ws.Range("T2").Value = 0.4626185066002986

# --- Row 3: MuSCs (new string, must be introduced first) | Gdf3 | Tdgf1 | (target set later) ---
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Gdf3"
$ws.Range("C3").Value = "Tdgf1"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.0822615
$ws.Range("H3").Value = 0.164523
$ws.Range("I3").Value = 0.002527782968212797
$ws.Range("J3").Value = 0.00168660977329485
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.08586166666666667
$ws.Range("N3").Value = 0.257585
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.007063109492500001
$ws.Range("R3").Value = 0.042378656955
$ws.Range("S3").Value = 0.002527782968212797
$ws.Range("T3").Value = 0.00168660977329485

# --- Row 4 (new row): Resolving-Mac | Gdf3 | Tdgf1 | (target set later) ---
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Gdf3"
$ws.Range("C4").Value = "Tdgf1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.41839966666667
$ws.Range("H4").Value = 52.255199
$ws.Range("I4").Value = 0.5352435101587433
$ws.Range("J4").Value = 0.5356948836264065
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.08586166666666667
$ws.Range("N4").Value = 0.257585
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1.495572826046111
$ws.Range("R4").Value = 13.460155434415
$ws.Range("S4").Value = 0.5352435101587433
$ws.Range("T4").Value = 0.5356948836264065

# --- Target cluster column: "ECs" is a brand-new shared string and must be
#     the LAST new string introduced (index 25), so set these last. ---
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "ECs"
$ws.Range("D4").Value = "ECs"
